$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V contents between rows 26/27 and 28/29 (columns A:E, i.e. Indice/
# pais/torneio/temporada/data_partida, stay exactly as they were) ---
# Row 26
$ws.Range("F26").Value = 'Atalanta'
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 'Monza'
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 1.56
$ws.Range("K26").Value = '20/08/2023 09:02'
$ws.Range("L26").Value = 1.74
$ws.Range("M26").Value = '02/09/2023 20:44'
$ws.Range("N26").Value = 4.31
$ws.Range("O26").Value = '20/08/2023 09:02'
$ws.Range("P26").Value = 4.09
$ws.Range("Q26").Value = '02/09/2023 20:44'
$ws.Range("R26").Value = 6.1
$ws.Range("S26").Value = '20/08/2023 09:02'
$ws.Range("T26").Value = 4.84
$ws.Range("U26").Value = '02/09/2023 20:44'
$ws.Range("V26").Value = 'https://www.betexplorer.com/football/italy/serie-a/atalanta-monza/4AhPuywD/'

# Row 27
$ws.Range("F27").Value = 'Napoli'
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 'Lazio'
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 1.72
$ws.Range("K27").Value = '20/08/2023 09:02'
$ws.Range("L27").Value = 1.67
$ws.Range("M27").Value = '02/09/2023 20:42'
$ws.Range("N27").Value = 3.87
$ws.Range("O27").Value = '20/08/2023 09:02'
$ws.Range("P27").Value = 4.18
$ws.Range("Q27").Value = '02/09/2023 20:44'
$ws.Range("R27").Value = 5.04
$ws.Range("S27").Value = '20/08/2023 09:02'
$ws.Range("T27").Value = 5.24
$ws.Range("U27").Value = '02/09/2023 20:43'
$ws.Range("V27").Value = 'https://www.betexplorer.com/football/italy/serie-a/napoli-lazio/vNZpcx8m/'

# Row 28
$ws.Range("F28").Value = 'Inter'
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 'Fiorentina'
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 1.72
$ws.Range("K28").Value = '20/08/2023 09:02'
$ws.Range("L28").Value = 1.74
$ws.Range("M28").Value = '03/09/2023 17:51'
$ws.Range("N28").Value = 3.94
$ws.Range("O28").Value = '20/08/2023 09:02'
$ws.Range("P28").Value = 3.99
$ws.Range("Q28").Value = '03/09/2023 18:09'
$ws.Range("R28").Value = 5
$ws.Range("S28").Value = '20/08/2023 09:02'
$ws.Range("T28").Value = 4.94
$ws.Range("U28").Value = '03/09/2023 18:09'
$ws.Range("V28").Value = 'https://www.betexplorer.com/football/italy/serie-a/inter-fiorentina/rilF2bpQ/'

# Row 29
$ws.Range("F29").Value = 'Torino'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 'Genoa'
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1.75
$ws.Range("K29").Value = '22/08/2023 14:46'
$ws.Range("L29").Value = 2.08
$ws.Range("M29").Value = '03/09/2023 18:29'
$ws.Range("N29").Value = 3.71
$ws.Range("O29").Value = '22/08/2023 14:46'
$ws.Range("P29").Value = 3.31
$ws.Range("Q29").Value = '03/09/2023 18:26'
$ws.Range("R29").Value = 4.59
$ws.Range("S29").Value = '22/08/2023 14:46'
$ws.Range("T29").Value = 4.16
$ws.Range("U29").Value = '03/09/2023 18:29'
$ws.Range("V29").Value = 'https://www.betexplorer.com/football/italy/serie-a/torino-genoa/UBScfzh6/'

# --- Append 6 new match rows (53-58) at the bottom of the table ---
# Row 53
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 'italy'
$ws.Range("C53").Value = 'serie-a'
$ws.Range("D53").Value = '2023-2024'
$ws.Range("E53").Value = 45196.77083333334
$ws.Range("F53").Value = 'Cagliari'
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 'AC Milan'
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = 5.55
$ws.Range("K53").Value = '23/09/2023 09:29'
$ws.Range("L53").Value = 4.73
$ws.Range("M53").Value = '27/09/2023 18:10'
$ws.Range("N53").Value = 4.14
$ws.Range("O53").Value = '23/09/2023 09:29'
$ws.Range("P53").Value = 3.43
$ws.Range("Q53").Value = '27/09/2023 18:10'
$ws.Range("R53").Value = 1.56
$ws.Range("S53").Value = '23/09/2023 09:29'
$ws.Range("T53").Value = 1.92
$ws.Range("U53").Value = '27/09/2023 18:10'
$ws.Range("V53").Value = 'https://www.betexplorer.com/football/italy/serie-a/cagliari-ac-milan/CWoooiWk/'
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null
$ws.Range("E52").Copy() | Out-Null
$ws.Range("E53").PasteSpecial(-4122) | Out-Null

# Row 54
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 'italy'
$ws.Range("C54").Value = 'serie-a'
$ws.Range("D54").Value = '2023-2024'
$ws.Range("E54").Value = 45196.77083333334
$ws.Range("F54").Value = 'Empoli'
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 'Salernitana'
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 2.25
$ws.Range("K54").Value = '17/09/2023 09:02'
$ws.Range("L54").Value = 2.8
$ws.Range("M54").Value = '27/09/2023 18:28'
$ws.Range("N54").Value = 3.02
$ws.Range("O54").Value = '17/09/2023 09:02'
$ws.Range("P54").Value = 3.19
$ws.Range("Q54").Value = '27/09/2023 18:26'
$ws.Range("R54").Value = 3.58
$ws.Range("S54").Value = '17/09/2023 09:02'
$ws.Range("T54").Value = 2.85
$ws.Range("U54").Value = '27/09/2023 18:28'
$ws.Range("V54").Value = 'https://www.betexplorer.com/football/italy/serie-a/empoli-salernitana/jNpkpBod/'
$ws.Range("A53").Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$ws.Range("E53").Copy() | Out-Null
$ws.Range("E54").PasteSpecial(-4122) | Out-Null

# Row 55
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 'italy'
$ws.Range("C55").Value = 'serie-a'
$ws.Range("D55").Value = '2023-2024'
$ws.Range("E55").Value = 45196.77083333334
$ws.Range("F55").Value = 'Verona'
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 'Atalanta'
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 3.96
$ws.Range("K55").Value = '17/09/2023 09:02'
$ws.Range("L55").Value = 4.48
$ws.Range("M55").Value = '27/09/2023 18:29'
$ws.Range("N55").Value = 3.76
$ws.Range("O55").Value = '17/09/2023 09:02'
$ws.Range("P55").Value = 3.66
$ws.Range("Q55").Value = '27/09/2023 18:29'
$ws.Range("R55").Value = 1.85
$ws.Range("S55").Value = '17/09/2023 09:02'
$ws.Range("T55").Value = 1.9
$ws.Range("U55").Value = '27/09/2023 18:29'
$ws.Range("V55").Value = 'https://www.betexplorer.com/football/italy/serie-a/verona-atalanta/4QAObA0k/'
$ws.Range("A54").Copy() | Out-Null
$ws.Range("A55").PasteSpecial(-4122) | Out-Null
$ws.Range("E54").Copy() | Out-Null
$ws.Range("E55").PasteSpecial(-4122) | Out-Null

# Row 56
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 'italy'
$ws.Range("C56").Value = 'serie-a'
$ws.Range("D56").Value = '2023-2024'
$ws.Range("E56").Value = 45196.86458333334
$ws.Range("F56").Value = 'Napoli'
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 'Udinese'
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = 1.35
$ws.Range("K56").Value = '17/09/2023 09:02'
$ws.Range("L56").Value = 1.43
$ws.Range("M56").Value = '27/09/2023 20:44'
$ws.Range("N56").Value = 5.18
$ws.Range("O56").Value = '17/09/2023 09:02'
$ws.Range("P56").Value = 5.07
$ws.Range("Q56").Value = '27/09/2023 20:43'
$ws.Range("R56").Value = 7.71
$ws.Range("S56").Value = '17/09/2023 09:02'
$ws.Range("T56").Value = 7.42
$ws.Range("U56").Value = '27/09/2023 20:44'
$ws.Range("V56").Value = 'https://www.betexplorer.com/football/italy/serie-a/napoli-udinese/hj7Kajoq/'
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null
$ws.Range("E55").Copy() | Out-Null
$ws.Range("E56").PasteSpecial(-4122) | Out-Null

# Row 57
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 'italy'
$ws.Range("C57").Value = 'serie-a'
$ws.Range("D57").Value = '2023-2024'
$ws.Range("E57").Value = 45196.86458333334
$ws.Range("F57").Value = 'Inter'
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 'Sassuolo'
$ws.Range("I57").Value = 2
$ws.Range("J57").Value = 1.41
$ws.Range("K57").Value = '17/09/2023 09:02'
$ws.Range("L57").Value = 1.28
$ws.Range("M57").Value = '27/09/2023 20:42'
$ws.Range("N57").Value = 5.06
$ws.Range("O57").Value = '17/09/2023 09:02'
$ws.Range("P57").Value = 6.79
$ws.Range("Q57").Value = '27/09/2023 20:44'
$ws.Range("R57").Value = 6.32
$ws.Range("S57").Value = '17/09/2023 09:02'
$ws.Range("T57").Value = 9.97
$ws.Range("U57").Value = '27/09/2023 20:43'
$ws.Range("V57").Value = 'https://www.betexplorer.com/football/italy/serie-a/inter-sassuolo/0hSu7Yh9/'
$ws.Range("A56").Copy() | Out-Null
$ws.Range("A57").PasteSpecial(-4122) | Out-Null
$ws.Range("E56").Copy() | Out-Null
$ws.Range("E57").PasteSpecial(-4122) | Out-Null

# Row 58
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 'italy'
$ws.Range("C58").Value = 'serie-a'
$ws.Range("D58").Value = '2023-2024'
$ws.Range("E58").Value = 45196.86458333334
$ws.Range("F58").Value = 'Lazio'
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 'Torino'
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1.89
$ws.Range("K58").Value = '17/09/2023 09:02'
$ws.Range("L58").Value = 2.31
$ws.Range("M58").Value = '27/09/2023 20:44'
$ws.Range("N58").Value = 3.5
$ws.Range("O58").Value = '17/09/2023 09:02'
$ws.Range("P58").Value = 3.27
$ws.Range("Q58").Value = '27/09/2023 20:44'
$ws.Range("R58").Value = 4.5
$ws.Range("S58").Value = '17/09/2023 09:02'
$ws.Range("T58").Value = 3.54
$ws.Range("U58").Value = '27/09/2023 20:44'
$ws.Range("V58").Value = 'https://www.betexplorer.com/football/italy/serie-a/lazio-torino/8IVm5CNL/'
$ws.Range("A57").Copy() | Out-Null
$ws.Range("A58").PasteSpecial(-4122) | Out-Null
$ws.Range("E57").Copy() | Out-Null
$ws.Range("E58").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
